$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1814.5588
$ws.Range("I112").Value = 1517
$ws.Range("J112").Value = 1865.862
$ws.Range("K112").Value = 4551
$ws.Range("L112").Value = 5597.586
$ws.Range("M112").Value = -3443
$ws.Range("N112").Value = -7813.586

$ws.Range("H113").Value = 7783.6313
$ws.Range("I113").Value = 6700.5713
$ws.Range("J113").Value = 8415.416999999999
$ws.Range("K113").Value = 6700.5713
$ws.Range("L113").Value = 8415.416999999999
$ws.Range("M113").Value = -3446.5713
$ws.Range("N113").Value = -14923.417

$ws.Range("H137").Value = 1893
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 1893
$ws.Range("K137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("M137").Value = 5679
$ws.Range("N137").Value = -10779

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 225
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 87.5
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 87.5
$ws.Range("M5").Value = -388
$ws.Range("N5").Value = -311.5

$ws.Range("H32").Value = 2799.8
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2799.8
$ws.Range("K32").Value = 0
$ws.Range("L32").ClearContents()
$ws.Range("M32").Value = 2799.8
$ws.Range("N32").Value = -3373.8

$ws.Range("H41").Value = 611.6667
$ws.Range("I41").Value = 517.5
$ws.Range("J41").Value = 800
$ws.Range("K41").Value = 517.5
$ws.Range("L41").Value = 800
$ws.Range("M41").Value = -103.5
$ws.Range("N41").Value = -1628

$ws.Range("H74").Value = 940
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 940
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("M74").Value = 940
$ws.Range("N74").Value = -2688

$ws.Range("H77").Value = 940
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 940
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("M77").Value = 4700
$ws.Range("N77").Value = -13436

$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0

$ws.Range("H122").Value = 1612.7142
$ws.Range("I122").Value = 1612.7142
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4838.142599999999
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -2388.142599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 225
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 87.5
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 87.5
$ws.Range("M4").Value = -385
$ws.Range("N4").Value = -317.5

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").ClearContents()
$ws.Range("N54").Value = 0

$ws.Range("H76").Value = 30003.5
$ws.Range("I76").Value = 25000
$ws.Range("J76").Value = 31671.334
$ws.Range("K76").Value = 25000
$ws.Range("L76").Value = 31671.334
$ws.Range("M76").Value = -24685
$ws.Range("N76").Value = -32301.334

$ws.Range("H79").Value = 30003.5
$ws.Range("I79").Value = 25000
$ws.Range("J79").Value = 31671.334
$ws.Range("K79").Value = 25000
$ws.Range("L79").Value = 31671.334
$ws.Range("M79").Value = -23908
$ws.Range("N79").Value = -33855.334

$ws.Range("H86").Value = 3317.0466
$ws.Range("I86").Value = 2645.625
$ws.Range("J86").Value = 5270.273
$ws.Range("K86").Value = 2645.625
$ws.Range("L86").Value = 5270.273
$ws.Range("M86").Value = -1522.625
$ws.Range("N86").Value = -7516.273

$ws.Range("H89").Value = 3317.0466
$ws.Range("I89").Value = 2645.625
$ws.Range("J89").Value = 5270.273
$ws.Range("K89").Value = 13228.125
$ws.Range("L89").Value = 26351.365
$ws.Range("M89").Value = -7612.125
$ws.Range("N89").Value = -37583.36500000001

$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0

$ws.Range("H107").Value = 12463.765
$ws.Range("I107").Value = 9777.429
$ws.Range("J107").Value = 25000
$ws.Range("K107").Value = 9777.429
$ws.Range("L107").Value = 25000
$ws.Range("M107").Value = -7857.429
$ws.Range("N107").Value = -28840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").ClearContents()
$ws.Range("N9").Value = 0

$ws.Range("H31").Value = 1837.2307
$ws.Range("I31").Value = 1679.1875
$ws.Range("J31").Value = 2559.7144
$ws.Range("K31").Value = 1679.1875
$ws.Range("L31").Value = 2559.7144
$ws.Range("M31").Value = -1384.1875
$ws.Range("N31").Value = -3149.7144

$ws.Range("H34").Value = 1837.2307
$ws.Range("I34").Value = 1679.1875
$ws.Range("J34").Value = 2559.7144
$ws.Range("K34").Value = 1679.1875
$ws.Range("L34").Value = 2559.7144
$ws.Range("M34").Value = -1477.1875
$ws.Range("N34").Value = -2963.7144

$ws.Range("H99").Value = 24937620
$ws.Range("I99").Value = 8131252
$ws.Range("J99").Value = 33340804
$ws.Range("K99").Value = 8131252
$ws.Range("L99").Value = 33340804
$ws.Range("M99").Value = -8129754
$ws.Range("N99").Value = -33343800

$ws.Range("H122").Value = 514678.4
$ws.Range("I122").Value = 852464
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 2557392
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -2554942
$ws.Range("N122").Value = -28900

$ws.Range("H126").Value = 24937620
$ws.Range("I126").Value = 8131252
$ws.Range("J126").Value = 33340804
$ws.Range("K126").Value = 24393756
$ws.Range("L126").Value = 100022412
$ws.Range("M126").Value = -24391286
$ws.Range("N126").Value = -100027352

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 96.71429000000001
$ws.Range("I2").Value = 98.59999999999999
$ws.Range("J2").Value = 92
$ws.Range("K2").Value = 591.5999999999999
$ws.Range("L2").Value = 552
$ws.Range("M2").Value = -478.5999999999999
$ws.Range("N2").Value = -778

$ws.Range("H34").Value = 699.8
$ws.Range("I34").Value = 99.666664
$ws.Range("J34").Value = 1600
$ws.Range("K34").Value = 298.999992
$ws.Range("L34").Value = 4800
$ws.Range("M34").Value = -214.999992
$ws.Range("N34").Value = -4968

$ws.Range("H101").Value = 3733.3333
$ws.Range("I101").Value = 1200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 800
$ws.Range("I6").Value = 800
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 800
$ws.Range("L6").Value = 800
$ws.Range("M6").Value = -687
$ws.Range("N6").Value = -1026

$ws.Range("H16").Value = 800
$ws.Range("I16").Value = 800
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = -550
$ws.Range("N16").Value = -1300

$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H122").Value = 1642.8572
$ws.Range("I122").Value = 1716.6666
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 5149.9998
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -2699.9998
$ws.Range("N122").Value = -8500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5660
$ws.Range("I93").Value = 6325
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 6325
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -5077
$ws.Range("N93").Value = -5496

$ws.Range("H116").Value = 236793.8
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 236793.8
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 236793.8
$ws.Range("N116").Value = -245971.8

$ws.Range("H122").Value = 5570.7915
$ws.Range("I122").Value = 2922.4
$ws.Range("J122").Value = 6267.737
$ws.Range("K122").Value = 8767.200000000001
$ws.Range("L122").Value = 18803.211
$ws.Range("M122").Value = -6317.200000000001
$ws.Range("N122").Value = -23703.211
